# Generate Report for Handoff
# adf73183-...md has been handed off for translation: its status moves from
# "In Translation" to "Ready for handoff" (priority -> "mt", new handoff
# datetime + handoff-file record). The regenerated report re-emits the rows
# so the still-in-translation file (f0ebc900-...md) now occupies row 2 and
# the just-handed-off file (adf73183-...md) occupies row 3 on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "f0ebc900-7944-4366-b695-87470ab2a147.md"
$ov.Range("A3").Value = "adf73183-2b88-4633-bbb6-f2b84e91e9b1.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-13 16:20:17"

$ov.Columns.Item(5).ColumnWidth = 17.2159881591797
$ov.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "f0ebc900-7944-4366-b695-87470ab2a147.md"
$zh.Range("G2").Value = "f0ebc900-7944-4366-b695-87470ab2a147.7d8044ddee7ddd5ab57723e1a903941f6dbd9835.zh-cn.xlf"

$zh.Range("A3").Value = "adf73183-2b88-4633-bbb6-f2b84e91e9b1.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("G3").Value = "adf73183-2b88-4633-bbb6-f2b84e91e9b1.cf793d2aa9b4de2547c2c3ab96d69d12d2788e9f.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-13 16:20:09"

$zh.Columns.Item(3).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "f0ebc900-7944-4366-b695-87470ab2a147.md"
$de.Range("G2").Value = "f0ebc900-7944-4366-b695-87470ab2a147.7d8044ddee7ddd5ab57723e1a903941f6dbd9835.de-de.xlf"

$de.Range("A3").Value = "adf73183-2b88-4633-bbb6-f2b84e91e9b1.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("G3").Value = "adf73183-2b88-4633-bbb6-f2b84e91e9b1.cf793d2aa9b4de2547c2c3ab96d69d12d2788e9f.de-de.xlf"
$de.Range("H3").Value = "2016-08-13 16:20:17"

$de.Columns.Item(3).ColumnWidth = 17.2159881591797
